$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for column F (dSF) following a repull/recalculation of the data.
$updates = @{
    15 = -4
    17 = 1
    18 = 0
    24 = -5
    25 = -9
    28 = 2
    29 = 2
    30 = -2
    32 = -8
    35 = -5
    41 = -2
    42 = 2
    46 = -11
    47 = 3
    51 = -4
    56 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
